$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 322
$ws.Range("I5").Value = 96.333336
$ws.Range("K5").Value = 96.333336
$ws.Range("M5").Value = 18.666664

$ws.Range("H17").Value = 1037.7646
$ws.Range("J17").Value = 1037.7646
$ws.Range("L17").Value = 3113.2938
$ws.Range("N17").Value = -3449.2938

$ws.Range("H19").Value = 1182.6666
$ws.Range("J19").Value = 1183
$ws.Range("L19").Value = 1183
$ws.Range("N19").Value = -1533

$ws.Range("H80").Value = 1856.4
$ws.Range("I80").Value = 2519.4
$ws.Range("J80").Value = 1524.9
$ws.Range("K80").Value = 7558.200000000001
$ws.Range("L80").Value = 4574.700000000001
$ws.Range("M80").Value = -6560.200000000001
$ws.Range("N80").Value = -6570.700000000001

$ws.Range("H83").Value = 1856.4
$ws.Range("I83").Value = 2519.4
$ws.Range("J83").Value = 1524.9
$ws.Range("K83").Value = 22674.6
$ws.Range("L83").Value = 13724.1
$ws.Range("M83").Value = -17682.6
$ws.Range("N83").Value = -23708.1

$ws.Range("H98").Value = 1473.5333
$ws.Range("I98").Value = 501.6
$ws.Range("J98").Value = 3417.4
$ws.Range("K98").Value = 501.6
$ws.Range("L98").Value = 3417.4
$ws.Range("M98").Value = 996.4
$ws.Range("N98").Value = -6413.4

$ws.Range("H107").Value = 1294.6428
$ws.Range("I107").Value = 1248.4
$ws.Range("K107").Value = 1248.4
$ws.Range("M107").Value = 671.5999999999999

$ws.Range("H122").Value = 1473.5333
$ws.Range("I122").Value = 501.6
$ws.Range("J122").Value = 3417.4
$ws.Range("K122").Value = 1504.8
$ws.Range("L122").Value = 10252.2
$ws.Range("M122").Value = 945.1999999999998
$ws.Range("N122").Value = -15152.2

$ws.Range("H132").Value = 4722.8
$ws.Range("I132").Value = 4708.2104
$ws.Range("K132").Value = 14124.6312
$ws.Range("M132").Value = -11594.6312

$ws.Range("H138").Value = 3278.4
$ws.Range("I138").Value = 698.5
$ws.Range("J138").Value = 4998.3335
$ws.Range("K138").Value = 2095.5
$ws.Range("L138").Value = 14995.0005
$ws.Range("M138").Value = 3044.5
$ws.Range("N138").Value = -25275.0005

$ws.Range("H141").Value = 10460.25
$ws.Range("I141").Value = 9934.091
$ws.Range("K141").Value = 29802.273
$ws.Range("M141").Value = -24622.273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 250.5
$ws.Range("I5").Value = 95
$ws.Range("K5").Value = 95
$ws.Range("M5").Value = 17

$ws.Range("H32").Value = 3894.1667
$ws.Range("I32").Value = 2262.1875
$ws.Range("J32").Value = 16950
$ws.Range("K32").Value = 2262.1875
$ws.Range("L32").Value = 16950
$ws.Range("M32").Value = -1975.1875
$ws.Range("N32").Value = -17524

$ws.Range("H74").Value = 1709.0769
$ws.Range("I74").Value = 1521.9
$ws.Range("K74").Value = 1521.9
$ws.Range("M74").Value = -647.9000000000001

$ws.Range("H77").Value = 1709.0769
$ws.Range("I77").Value = 1521.9
$ws.Range("K77").Value = 7609.5
$ws.Range("M77").Value = -3241.5

$ws.Range("H102").Value = 35716010
$ws.Range("I102").Value = 35716010
$ws.Range("K102").Value = 35716010
$ws.Range("M102").Value = -35714388

$ws.Range("H122").Value = 1379.2307
$ws.Range("I122").Value = 1448.4572
$ws.Range("J122").Value = 773.5
$ws.Range("K122").Value = 4345.3716
$ws.Range("L122").Value = 2320.5
$ws.Range("M122").Value = -1895.3716
$ws.Range("N122").Value = -7220.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 250.5
$ws.Range("I4").Value = 95
$ws.Range("K4").Value = 95
$ws.Range("M4").Value = 20

$ws.Range("H20").Value = 1485.7142
$ws.Range("K20").Value = 1600
$ws.Range("M20").Value = -1353

$ws.Range("H107").Value = 3268.4614
$ws.Range("I107").Value = 3268.4614
$ws.Range("K107").Value = 3268.4614
$ws.Range("M107").Value = -1348.4614

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 2833
$ws.Range("J29").Value = 2749.5
$ws.Range("L29").Value = 2749.5
$ws.Range("N29").Value = -3335.5

$ws.Range("H58").Value = 3063.9
$ws.Range("I58").Value = 1806
$ws.Range("K58").Value = 1806
$ws.Range("M58").Value = -1603

$ws.Range("H122").Value = 1532.5333
$ws.Range("I122").Value = 1532.5333
$ws.Range("K122").Value = 4597.5999
$ws.Range("M122").Value = -2147.5999

$ws.Range("H132").Value = 2458.8125
$ws.Range("I132").Value = 2289.4666
$ws.Range("K132").Value = 6868.399800000001
$ws.Range("M132").Value = -4338.399800000001

$ws.Range("H134").Value = 3400.7856
$ws.Range("I134").Value = 3440.8462
$ws.Range("J134").Value = 2880
$ws.Range("K134").Value = 10322.5386
$ws.Range("L134").Value = 8640
$ws.Range("M134").Value = -7787.5386
$ws.Range("N134").Value = -13710

$ws.Range("H136").Value = 3063.9
$ws.Range("I136").Value = 1806
$ws.Range("K136").Value = 5418
$ws.Range("M136").Value = -2868

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4782556
$ws.Range("I4").Value = 2842752.5
$ws.Range("K4").Value = 8528257.5
$ws.Range("M4").Value = -8528145.5

$ws.Range("H47").Value = 402.85715
$ws.Range("I47").Value = 164.4
$ws.Range("J47").Value = 999
$ws.Range("K47").Value = 493.2
$ws.Range("L47").Value = 2997
$ws.Range("M47").Value = -62.20000000000005
$ws.Range("N47").Value = -3859

$ws.Range("H69").Value = 788
$ws.Range("I69").Value = 835
$ws.Range("K69").Value = 2505
$ws.Range("M69").Value = -1694

$ws.Range("H72").Value = 788
$ws.Range("I72").Value = 835
$ws.Range("K72").Value = 7515
$ws.Range("M72").Value = -3459

$ws.Range("H92").Value = 333.33334
$ws.Range("I92").Value = 350
$ws.Range("J92").Value = 325
$ws.Range("K92").Value = 1050
$ws.Range("L92").Value = 975
$ws.Range("M92").Value = 198
$ws.Range("N92").Value = -3471

$ws.Range("H107").Value = 514.8182
$ws.Range("J107").Value = 480
$ws.Range("L107").Value = 1440
$ws.Range("N107").Value = -5280

$ws.Range("H137").Value = 5343.0835
$ws.Range("I137").Value = 4370
$ws.Range("J137").Value = 5537.7
$ws.Range("K137").Value = 13110
$ws.Range("L137").Value = 16613.1
$ws.Range("M137").Value = -8010
$ws.Range("N137").Value = -26813.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3276.2
$ws.Range("I102").Value = 3484.6667
$ws.Range("J102").Value = 1400
$ws.Range("K102").Value = 3484.6667
$ws.Range("L102").Value = 1400
$ws.Range("M102").Value = -1862.6667
$ws.Range("N102").Value = -4644

$ws.Range("H132").Value = 2499.3333
$ws.Range("I132").Value = 2499.5
$ws.Range("J132").Value = 2499
$ws.Range("K132").Value = 7498.5
$ws.Range("L132").Value = 7497
$ws.Range("M132").Value = -4968.5
$ws.Range("N132").Value = -12557

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1378
$ws.Range("I16").Value = 823.3333
$ws.Range("K16").Value = 823.3333
$ws.Range("M16").Value = -653.3333

$ws.Range("H22").Value = 2622.111
$ws.Range("I22").Value = 1979.8
$ws.Range("K22").Value = 1979.8
$ws.Range("M22").Value = -1684.8

$ws.Range("H23").Value = 5000
$ws.Range("I23").Value = 5000
$ws.Range("K23").Value = 5000
$ws.Range("M23").Value = -4770

$ws.Range("H25").Value = 2000
$ws.Range("I25").Value = 2000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 2000
$ws.Range("M25").Value = -1770
$ws.Range("N25").Value = 0

$ws.Range("H27").Value = 2622.111
$ws.Range("I27").Value = 1979.8
$ws.Range("K27").Value = 1979.8
$ws.Range("M27").Value = -1872.8

$ws.Range("H40").Value = 4962.375
$ws.Range("I40").Value = 4915.3076
$ws.Range("K40").Value = 4915.3076
$ws.Range("M40").Value = -4779.3076

$ws.Range("H46").Value = 3725.889
$ws.Range("I46").Value = 3063
$ws.Range("J46").Value = 4057.3333
$ws.Range("K46").Value = 3063
$ws.Range("L46").Value = 4057.3333
$ws.Range("M46").Value = -2875
$ws.Range("N46").Value = -4433.3333

$ws.Range("H55").Value = 1226.3529
$ws.Range("I55").Value = 1240
$ws.Range("J55").Value = 1220.6666
$ws.Range("K55").Value = 1240
$ws.Range("L55").Value = 1220.6666
$ws.Range("M55").Value = -1067
$ws.Range("N55").Value = -1566.6666

$ws.Range("H93").Value = 2656.5715
$ws.Range("I93").Value = 3167
$ws.Range("J93").Value = 2273.75
$ws.Range("K93").Value = 3167
$ws.Range("L93").Value = 2273.75
$ws.Range("M93").Value = -1919
$ws.Range("N93").Value = -4769.75

$ws.Range("H136").Value = 17329.666
$ws.Range("I136").Value = 17244.5
$ws.Range("J136").Value = 17500
$ws.Range("K136").Value = 51733.5
$ws.Range("L136").Value = 52500
$ws.Range("M136").Value = -49183.5
$ws.Range("N136").Value = -57600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 45000
$ws.Range("J27").Value = 45000
$ws.Range("L27").Value = 45000
$ws.Range("N27").Value = -45138

$ws.Range("H32").Value = 3338138
$ws.Range("I32").Value = 4005379.8
$ws.Range("K32").Value = 4005379.8
$ws.Range("M32").Value = -4005062.8

$ws.Range("H45").Value = 34310.89
$ws.Range("J45").Value = 36977.5
$ws.Range("L45").Value = 36977.5
$ws.Range("N45").Value = -37959.5

$ws.Range("H54").Value = 38000
$ws.Range("J54").Value = 50000
$ws.Range("L54").Value = 50000
$ws.Range("N54").Value = -51040

$ws.Range("H126").Value = 2181.05
$ws.Range("I126").Value = 2263.9375
$ws.Range("K126").Value = 6791.8125
$ws.Range("M126").Value = -4321.8125

$ws.Range("H132").Value = 3184
$ws.Range("I132").Value = 2349.3572
$ws.Range("K132").Value = 7048.071599999999
$ws.Range("M132").Value = -4518.071599999999
